$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$newValues = @(
    "Stimuli/327.jpg",
    "Stimuli/1111.jpg",
    "Stimuli/3017.jpg",
    "Stimuli/3022.jpg",
    "Stimuli/3180.jpg",
    "Stimuli/3280.jpg",
    "Stimuli/6190.jpg",
    "Stimuli/6244.jpg",
    "Stimuli/6836.jpg",
    "Stimuli/9180.jpg",
    "Stimuli/9182.jpg",
    "Stimuli/9253.jpg",
    "Stimuli/9300.jpg",
    "Stimuli/9326.jpg",
    "Stimuli/9424.jpg",
    "Stimuli/9425.jpg",
    "Stimuli/9920.jpg"
)

$row = 5
foreach ($val in $newValues) {
    $ws.Cells.Item($row, 1).Value = $val
    $row = $row + 1
}

$ws.Rows.Item(22).Select()
